# "Generate Report for Handoff" -- refresh the localization-status report
# with the latest handoff run's identifiers/timestamps.
#
# Old run id: 80bb7ea5-eed1-4f96-bcf6-9a2d6aaa2c5d
# New run id: a31496fd-c88c-4a94-b6dd-2d76bf0f6dee
#
# Old xlf hash: 91e3b33e1bcac472ca4b010bc2efcf1d7a48f6af
# New xlf hash: b00050c3ad61f5a517547bfc4158803d6d679d57

$wb = $excel.ActiveWorkbook

$oldId = "80bb7ea5-eed1-4f96-bcf6-9a2d6aaa2c5d"
$newId = "a31496fd-c88c-4a94-b6dd-2d76bf0f6dee"

$oldMd = "$oldId.md"
$newMd = "$newId.md"

$oldXlfZh = "$oldId.91e3b33e1bcac472ca4b010bc2efcf1d7a48f6af.zh-cn.xlf"
$newXlfZh = "$newId.b00050c3ad61f5a517547bfc4158803d6d679d57.zh-cn.xlf"

$oldXlfDe = "$oldId.91e3b33e1bcac472ca4b010bc2efcf1d7a48f6af.de-de.xlf"
$newXlfDe = "$newId.b00050c3ad61f5a517547bfc4158803d6d679d57.de-de.xlf"

$oldOverviewDate = "2016-47-12 16:47:57"
$newOverviewDate = "2016-48-12 16:48:16"

$oldZhDate = "2016-03-12 16:47:54"
$newZhDate = "2016-03-12 16:48:12"

$oldDeDate = "2016-03-12 16:47:57"
$newDeDate = "2016-03-12 16:48:16"

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = $newOverviewDate
foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.TextToDisplay -eq $oldMd) {
        $h.TextToDisplay = $newMd
    }
}

# ---- Sheet "zh-cn" ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMd
$wsZhCn.Range("D2").Value = $newXlfZh
$wsZhCn.Range("E2").Value = $newZhDate
foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.TextToDisplay -eq $oldMd) {
        $h.TextToDisplay = $newMd
    } elseif ($h.TextToDisplay -eq $oldXlfZh) {
        $h.TextToDisplay = $newXlfZh
    }
}

# ---- Sheet "de-de" ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMd
$wsDeDe.Range("D2").Value = $newXlfDe
$wsDeDe.Range("E2").Value = $newDeDate
foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.TextToDisplay -eq $oldMd) {
        $h.TextToDisplay = $newMd
    } elseif ($h.TextToDisplay -eq $oldXlfDe) {
        $h.TextToDisplay = $newXlfDe
    }
}
